# Example Data 2 - Amelia Tristan
# Added two variables: inseam length in cm, and hair color
# (Replaces the previous "Eye color" / "Waist" variables with
#  "Inseam" / "Hair Color" on the Data sheet, and updates the
#  Codebook sheet to match.)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Data"
$ws2 = $wb.Worksheets.Item(2)   # "Codebook"

# --- Data sheet: rename column headers ----------------------------------
$ws1.Cells.Item(1,4).Value = "Inseam"
$ws1.Cells.Item(1,5).Value = "Hair Color"

# --- Data sheet: replace column D (Eye color -> Inseam, numeric cm) -----
$ws1.Cells.Item(2,4).Value = 81
$ws1.Cells.Item(3,4).Value = 74
$ws1.Cells.Item(4,4).Value = 55
$ws1.Cells.Item(5,4).Value = 91
$ws1.Cells.Item(6,4).Value = 95
$ws1.Cells.Item(7,4).Value = 89
$ws1.Cells.Item(8,4).Value = 68
$ws1.Cells.Item(9,4).Value = 62
$ws1.Cells.Item(10,4).Value = 73
$ws1.Cells.Item(11,4).Value = 49
$ws1.Cells.Item(12,4).Value = 53
$ws1.Cells.Item(13,4).Value = 50
$ws1.Cells.Item(14,4).Value = 52
$ws1.Cells.Item(15,4).Value = 65

# --- Data sheet: replace column E (Waist -> Hair Color, text) -----------
$ws1.Cells.Item(2,5).Value = "black"
$ws1.Cells.Item(3,5).Value = "blond"
$ws1.Cells.Item(4,5).Value = "brown"
$ws1.Cells.Item(5,5).Value = "l brown"
$ws1.Cells.Item(6,5).Value = "other"
$ws1.Cells.Item(7,5).Value = "black"
$ws1.Cells.Item(8,5).Value = "other"
$ws1.Cells.Item(9,5).Value = "blond"
$ws1.Cells.Item(10,5).Value = "brown"
$ws1.Cells.Item(11,5).Value = "other"
$ws1.Cells.Item(12,5).Value = "other"
$ws1.Cells.Item(13,5).Value = "l brown"
$ws1.Cells.Item(14,5).Value = "d brown"
$ws1.Cells.Item(15,5).Value = "white"

# --- Codebook sheet: drop the now-removed "Eye Color" / "Waist" rows ----
# (Row 6 first so row numbers of remaining rows don't shift out from
#  under us before the second delete.)
$ws2.Rows.Item(6).Delete()
$ws2.Rows.Item(5).Delete()

# --- Selection / active sheet state, matching the saved workbook --------
$ws2.Range("A4").Select()
$ws1.Activate()
$ws1.Range("E15").Select()
